$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row5 = $ws.Rows.Item(5)
$row5.RowHeight = 92.25

$r5 = $ws.Range("D5")
$r5.Font.Underline = $true
$r5.VerticalAlignment = -4160
$r5.WrapText = $true

$ws.Range("A5:XFD5").Select()

Write-Output "ok"
